$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 725
$ws.Range("I18").Value = 725
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 725
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -441
$ws.Range("N18").ClearContents()
$ws.Range("H132").Value = 15916.823
$ws.Range("I132").Value = 19444.182
$ws.Range("J132").Value = 9450
$ws.Range("K132").Value = 58332.546
$ws.Range("L132").Value = 28350
$ws.Range("M132").Value = -55802.546
$ws.Range("N132").Value = -33410
$ws.Range("H137").Value = 2848.762
$ws.Range("I137").Value = 2011.909
$ws.Range("K137").Value = 6035.727000000001
$ws.Range("M137").Value = -3485.727000000001
$ws.Range("H138").Value = 2059.1
$ws.Range("I138").Value = 1755.875
$ws.Range("J138").Value = 3272
$ws.Range("K138").Value = 5267.625
$ws.Range("L138").Value = 9816
$ws.Range("M138").Value = -127.625
$ws.Range("N138").Value = -20096

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3567.8572
$ws.Range("I45").Value = 2795
$ws.Range("K45").Value = 2795
$ws.Range("M45").Value = -2418
$ws.Range("H74").Value = 7242.6665
$ws.Range("J74").Value = 9921.5
$ws.Range("L74").Value = 9921.5
$ws.Range("N74").Value = -11669.5
$ws.Range("H77").Value = 7242.6665
$ws.Range("J77").Value = 9921.5
$ws.Range("L77").Value = 49607.5
$ws.Range("N77").Value = -58343.5
$ws.Range("H122").Value = 1994
$ws.Range("I122").Value = 1994
$ws.Range("K122").Value = 5982
$ws.Range("M122").Value = -3532
$ws.Range("H132").Value = 2587.2144
$ws.Range("J132").Value = 4569.8
$ws.Range("L132").Value = 13709.4
$ws.Range("N132").Value = -18769.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H29").Value = 688.6667
$ws.Range("I29").Value = 688.6667
$ws.Range("K29").Value = 688.6667
$ws.Range("M29").Value = -399.6667
$ws.Range("H99").Value = 4835.375
$ws.Range("I99").Value = 4540.8
$ws.Range("K99").Value = 4540.8
$ws.Range("M99").Value = -3042.8
$ws.Range("H134").Value = 3612.35
$ws.Range("I134").Value = 2883.2
$ws.Range("K134").Value = 8649.599999999999
$ws.Range("M134").Value = -6114.599999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5883.6665
$ws.Range("I31").Value = 1231.4445
$ws.Range("K31").Value = 1231.4445
$ws.Range("M31").Value = -936.4445000000001
$ws.Range("H34").Value = 5883.6665
$ws.Range("I34").Value = 1231.4445
$ws.Range("K34").Value = 1231.4445
$ws.Range("M34").Value = -1029.4445
$ws.Range("H86").Value = 2002.3334
$ws.Range("I86").Value = 2003.5
$ws.Range("K86").Value = 2003.5
$ws.Range("M86").Value = -880.5
$ws.Range("H89").Value = 2002.3334
$ws.Range("I89").Value = 2003.5
$ws.Range("K89").Value = 10017.5
$ws.Range("M89").Value = -4401.5
$ws.Range("H132").Value = 4248
$ws.Range("I132").Value = 3383
$ws.Range("J132").Value = 7996.3335
$ws.Range("K132").Value = 10149
$ws.Range("L132").Value = 23989.0005
$ws.Range("M132").Value = -7619
$ws.Range("N132").Value = -29049.0005
$ws.Range("H134").Value = 1000
$ws.Range("I134").Value = 1000
$ws.Range("K134").Value = 3000
$ws.Range("M134").Value = -465

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 373.625
$ws.Range("I46").Value = 279.8
$ws.Range("J46").Value = 530
$ws.Range("K46").Value = 839.4000000000001
$ws.Range("L46").Value = 1590
$ws.Range("M46").Value = -748.4000000000001
$ws.Range("N46").Value = -1772
$ws.Range("H81").Value = 1869.3334
$ws.Range("J81").Value = 1869.3334
$ws.Range("L81").Value = 5608.0002
$ws.Range("N81").Value = -7854.0002
$ws.Range("H84").Value = 1869.3334
$ws.Range("J84").Value = 1869.3334
$ws.Range("L84").Value = 16824.0006
$ws.Range("N84").Value = -28056.0006
$ws.Range("H98").Value = 154.4
$ws.Range("I98").Value = 154.4
$ws.Range("K98").Value = 463.2
$ws.Range("M98").Value = 1034.8
$ws.Range("H118").Value = 866.6667
$ws.Range("J118").Value = 2000
$ws.Range("L118").Value = 6000
$ws.Range("N118").Value = -8486
$ws.Range("H121").Value = 696.2222
$ws.Range("J121").Value = 1729
$ws.Range("L121").Value = 5187
$ws.Range("N121").Value = -7807
$ws.Range("H131").Value = 2148.4546
$ws.Range("J131").Value = 2132.6667
$ws.Range("L131").Value = 6398.000100000001
$ws.Range("N131").Value = -16478.0001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H55").Value = 5439.625
$ws.Range("I55").Value = 5821.3335
$ws.Range("J55").Value = 4294.5
$ws.Range("K55").Value = 5821.3335
$ws.Range("L55").Value = 4294.5
$ws.Range("M55").Value = -5494.3335
$ws.Range("N55").Value = -4948.5
$ws.Range("H132").Value = 98352.63
$ws.Range("I132").Value = 107237.9
$ws.Range("K132").Value = 321713.7
$ws.Range("M132").Value = -319183.7

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H23").Value = 70000000
$ws.Range("I23").Value = 70000000
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 70000000
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = -69999770
$ws.Range("N23").ClearContents()
$ws.Range("H25").Value = 7257.5
$ws.Range("H29").Value = 50000
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 50000
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 50000
$ws.Range("M29").ClearContents()
$ws.Range("N29").Value = -50590
$ws.Range("H132").Value = 5098.75
$ws.Range("I132").Value = 3935
$ws.Range("J132").Value = 8590
$ws.Range("K132").Value = 11805
$ws.Range("L132").Value = 25770
$ws.Range("M132").Value = -9275
$ws.Range("N132").Value = -30830

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 2817500
$ws.Range("I21").Value = 5000000
$ws.Range("J21").Value = 2090000
$ws.Range("K21").Value = 5000000
$ws.Range("L21").Value = 2090000
$ws.Range("M21").Value = -4999765
$ws.Range("N21").Value = -2090470
$ws.Range("H24").Value = 10000
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 10000
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 10000
$ws.Range("M24").ClearContents()
$ws.Range("N24").Value = -10460
$ws.Range("H29").Value = 5000
$ws.Range("I29").Value = 5000
$ws.Range("K29").Value = 5000
$ws.Range("M29").Value = -4710
$ws.Range("H30").Value = 33505
$ws.Range("J30").Value = 33505
$ws.Range("L30").Value = 33505
$ws.Range("N30").Value = -33719
$ws.Range("H35").Value = 2817500
$ws.Range("I35").Value = 5000000
$ws.Range("J35").Value = 2090000
$ws.Range("K35").Value = 5000000
$ws.Range("L35").Value = 2090000
$ws.Range("M35").Value = -4999710
$ws.Range("N35").Value = -2090580
$ws.Range("H64").Value = 60000
$ws.Range("J64").Value = 60000
$ws.Range("L64").Value = 60000
$ws.Range("N64").Value = -60496
$ws.Range("H67").Value = 60000
$ws.Range("J67").Value = 60000
$ws.Range("L67").Value = 60000
$ws.Range("N67").Value = -61716
$ws.Range("H132").Value = 2819.5
$ws.Range("I132").Value = 2507.4666
$ws.Range("K132").Value = 7522.399800000001
$ws.Range("M132").Value = -4992.399800000001
$ws.Range("H136").Value = 4323.5454
$ws.Range("I136").Value = 3833.3333
$ws.Range("J136").Value = 4911.8
$ws.Range("K136").Value = 11499.9999
$ws.Range("L136").Value = 14735.4
$ws.Range("M136").Value = -8949.999899999999
$ws.Range("N136").Value = -19835.4

